$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date (G2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 13:10:38"

# zh-cn sheet: Correspond Handoff Datetime (H2), Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 13:10:34"
$wsZhCn.Range("K2").Value = "2016-09-04 13:10:51"

# de-de sheet: Correspond Handoff Datetime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 13:10:58"
